$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 0.8598786354551747
$ws.Range("D2").Value = 0.3958786986172012

# Row 3
$ws.Range("C3").Value = -2.166286983853513
$ws.Range("D3").Value = 0.03739143186762073

# Row 4
$ws.Range("C4").Value = -1.394629964818566
$ws.Range("D4").Value = 0.1721734310753793

# Row 5
$ws.Range("C5").Value = -2.413184202646503
$ws.Range("D5").Value = 0.02135197105876241

# Row 6
$ws.Range("C6").Value = -2.0327767670436
$ws.Range("D6").Value = 0.04994332885681874

# Row 7
$ws.Range("C7").Value = -1.729559377264682
$ws.Range("D7").Value = 0.09277951885676061

# Row 8
$ws.Range("C8").Value = -2.273979684159579
$ws.Range("D8").Value = 0.02939808063263105

# Row 9
$ws.Range("C9").Value = 0.6526159095817147
$ws.Range("D9").Value = 0.5183924446284658

# Row 10
$ws.Range("C10").Value = -2.024487921746192
$ws.Range("D10").Value = 0.05083227895260101
$ws.Range("G10").Value = "No"

# Row 11
$ws.Range("C11").Value = -2.000417524789786
$ws.Range("D11").Value = 0.05349263671563342
$ws.Range("G11").Value = "No"
